# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the d86d483f-... file after a fresh handback report was generated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: refresh the "Latest HO Xliff Generate Date" column for row 2 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-21 14:55:20"

# --- zh-cn sheet: refresh handoff/handback datetimes for row 2 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-21 14:55:15"
$wsZhCn.Range("K2").Value = "2016-08-21 14:55:32"

# --- de-de sheet: refresh handoff/handback datetimes for row 2 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-21 14:55:20"
$wsDeDe.Range("K2").Value = "2016-08-21 14:55:38"
